# Added HTML and TXT handling
# Update the "cars2" sheet data (carModel / price / color) with new sample rows,
# and refresh the active-cell selections that Excel persisted on each sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # cars
$ws2 = $wb.Worksheets.Item(2)   # cars2
$ws3 = $wb.Worksheets.Item(3)   # cars3

# --- cars2 data changes ------------------------------------------------
$ws2.Range("A2").Value = "Jeep"
$ws2.Range("B2").Value = 25000
$ws2.Range("C2").Value = "purple"

$ws2.Range("A3").Value = "Voltswagon"
$ws2.Range("B3").Value = 16000
$ws2.Range("C3").Value = "black"

$ws2.Range("A4").Value = "Nissan"
$ws2.Range("B4").Value = 13000
$ws2.Range("C4").Value = "yellow"

$ws2.Range("A5").Value = "Honda"
$ws2.Range("B5").Value = 20000
$ws2.Range("C5").Value = "green"

# --- selection / active-tab bookkeeping --------------------------------
$ws2.Range("E11").Select() | Out-Null
$ws3.Range("F23").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("G34").Select() | Out-Null
